# Switches - front panel - wip
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 - the existing "J2" label is renamed to "J1" (first jack row keeps its
# original position values; a fresh "J2" row is inserted right after it).
$ws.Range("A7").Value = "J1"

# Row 8 - new "J2" entry
$ws.Range("A8").Value = "J2"
$ws.Range("B8").Value = 115.61750000000001
$ws.Range("C8").Value = 38.384
$ws.Range("D8").Formula = '=B8+$B$5'
$ws.Range("E8").Formula = '=C8+$C$5'

# Row 9 - "J3"
$ws.Range("A9").Value = "J3"
$ws.Range("B9").Value = 105.9
$ws.Range("C9").Value = 69.700666999999996
$ws.Range("D9").Formula = '=B9+$B$5'
$ws.Range("E9").Formula = '=C9+$C$5'

# Row 10 - "J4"
$ws.Range("A10").Value = "J4"
$ws.Range("D10").Formula = '=B10+$B$5'
$ws.Range("E10").Formula = '=C10+$C$5'

# Row 11 - "J5"
$ws.Range("A11").Value = "J5"
$ws.Range("C11").Value = 101.01733299999999
$ws.Range("D11").Formula = '=B11+$B$5'
$ws.Range("E11").Formula = '=C11+$C$5'

# Row 12 - "J6"
$ws.Range("A12").Value = "J6"
$ws.Range("D12").Formula = '=B12+$B$5'
$ws.Range("E12").Formula = '=C12+$C$5'

# Row 13 - "J7"
$ws.Range("A13").Value = "J7"
$ws.Range("C13").Value = 132.334
$ws.Range("D13").Formula = '=B13+$B$5'
$ws.Range("E13").Formula = '=C13+$C$5'

# Row 14 - "J8"
$ws.Range("A14").Value = "J8"
$ws.Range("D14").Formula = '=B14+$B$5'
$ws.Range("E14").Formula = '=C14+$C$5'

# Row 15 - "SW1"
$ws.Range("A15").Value = "SW1"
$ws.Range("B15").Value = 131.97999999999999
$ws.Range("C15").Value = 42.457999999999998
$ws.Range("D15").Formula = '=B15+$B$6'
$ws.Range("E15").Formula = '=C15+$C$6'

# Row 16 - "SW2"
$ws.Range("A16").Value = "SW2"
$ws.Range("C16").Value = 73.774666999999994
$ws.Range("D16").Formula = '=B16+$B$6'
$ws.Range("E16").Formula = '=C16+$C$6'

# Row 17 - "SW3"
$ws.Range("A17").Value = "SW3"
$ws.Range("C17").Value = 105.09133300000001
$ws.Range("D17").Formula = '=B17+$B$6'
$ws.Range("E17").Formula = '=C17+$C$6'

# Row 18 - "SW4"
$ws.Range("A18").Value = "SW4"
$ws.Range("C18").Value = 136.40799999999999
$ws.Range("D18").Formula = '=B18+$B$6'
$ws.Range("E18").Formula = '=C18+$C$6'

# Match the final selection recorded in the saved workbook
$ws.Range("E18").Select()
